$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 211, shifting existing rows 211-317 down to 212-318.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with its data (mirrors the row that used
# to be 211, now at 212, but with updated Fecha/Volumen/Precios/Precio $/Kg).
$ws.Range("A211").Value = 4
$ws.Range("B211").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C211").Value = "Los Lagos"
$ws.Range("D211").Value = 44806
$ws.Range("E211").Value = 10
$ws.Range("F211").Value = 100112043
$ws.Range("G211").Value = "Pepino ensalada"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 400
$ws.Range("K211").Value = 30000
$ws.Range("L211").Value = 30000
$ws.Range("M211").Value = 30000
$ws.Range("N211").Value = "`$/caja 60 unidades"
$ws.Range("O211").Value = "Región de Arica y Parinacota"
$ws.Range("P211").Value = 500
$ws.Range("Q211").Value = 60
$ws.Range("R211").Value = "Hortaliza"
